# Fixed some errors in the Gaussian weighted fits - more clean
#
# Row 8 corresponds to cellID "22-05-31_M7_P1_C8":
#   - "excluded" (column D) changes from "yes" to "no"
#   - "comments" (column E) is cleared (was "shitty curve")
# Also move the active cell selection to E8 (reflecting where the edit was made).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update excluded flag for row 8
$ws.Range("D8").Value = "no"

# Clear the comment in E8 entirely (not just blank text, remove the cell's content)
$ws.Range("E8").ClearContents()

# Update the selection shown in the sheet view to E8
$ws.Range("E8").Select()
